$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (rows 2-51) to text format so that numeric-looking
# strings (e.g. "1.000", "0.3243") are preserved as text instead of being
# parsed into numbers, matching the original inline-string cell type.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.061.31'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '1.912.21'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '0.8296'
$ws.Range("E5").Value = '  +8.46%  '
$ws.Range("D6").Value = '242.50'
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '0.3243'
$ws.Range("E8").Value = '  +5.57%  '
$ws.Range("D9").Value = '26.80'
$ws.Range("E9").Value = '  +4.25%  '
$ws.Range("D10").Value = '0.07042'
$ws.Range("E10").Value = '  +2.77%  '
$ws.Range("D11").Value = '0.08044'
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("D12").Value = '0.7532'
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").Value = '1.905.30'
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").Value = '5.239'
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").Value = '92.97'
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").Value = '14.22'
$ws.Range("E16").Value = '  +1.62%  '
$ws.Range("D17").Value = '30.060.83'
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("D18").Value = '5.961'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").Value = '245.24'
$ws.Range("E19").Value = '  +0.83%  '
$ws.Range("D20").Value = '0.000007785'
$ws.Range("E20").Value = '  +1.05%  '
$ws.Range("D21").Value = '2.156.97'
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '7.011'
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").Value = '0.1639'
$ws.Range("E25").Value = '  +26.27%  '
$ws.Range("D26").Value = '169.98'
$ws.Range("E26").Value = '  +1.96%  '
$ws.Range("D27").Value = '9.281'
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").Value = '18.98'
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("D29").Value = '2.089'
$ws.Range("E29").Value = '  +2.22%  '
$ws.Range("D30").Value = '1.373'
$ws.Range("E30").Value = '  -2.43%  '
$ws.Range("D31").Value = '1.521'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").Value = '4.312'
$ws.Range("D33").Value = '0.05619'
$ws.Range("E33").Value = '  +6.09%  '
$ws.Range("D34").Value = '4.105'
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("D35").Value = '1.287'
$ws.Range("E35").Value = '  +2.49%  '
$ws.Range("D36").Value = '0.7372'
$ws.Range("E36").Value = '  +1.15%  '
$ws.Range("D37").Value = '2.717'
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").Value = '0.01923'
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("D39").Value = '2.793'
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("D40").Value = '0.4450'
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("D41").Value = '72.71'
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("D42").Value = '6.011'
$ws.Range("E42").Value = '  -2.84%  '
$ws.Range("D43").Value = '0.8436'
$ws.Range("E43").Value = '  +1.54%  '
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '1.905'
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("D46").Value = '7.633'
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").Value = '101.34'
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("D48").Value = '9.765'
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("D49").Value = '987.07'
$ws.Range("E49").Value = '  +9.33%  '
$ws.Range("D50").Value = '2.063.86'
$ws.Range("E50").Value = '  +0.32%  '
$ws.Range("D51").Value = '36.38'
$ws.Range("E51").Value = '  +0.59%  '

# Restore the default "Normal" style so no extra number-format style is
# left behind on these cells.
$dataRange.Style = "Normal"
